# Nexial "number-showcase.xlsx" maintenance edit
# -------------------------------------------------------------
# The hidden "#system" sheet holds, per-column, an alphabetically
# sorted catalogue of command names for each command group (the
# defined names "base", "external", "io", "web", ... each point at
# one such column). This change set:
#   - removes "clear(variables)" from the "base" catalogue (col F)
#   - adds "terminate(programName)" to the "external" catalogue (col J)
#   - adds "assertPath(path)" to the "io" catalogue (col L)
#   - corrects "assertAttributeContains" -> "assertAttributeContain"
#     and adds "saveSelectedText"/"saveSelectedValue" to the "web"
#     catalogue (col Z)
# and updates the defined-name ranges + used-range dimension that
# describe those columns to reflect the new row counts.
#
# NOTE: this runtime's Range.Insert()/Range.Delete() shift the WHOLE
# row (every column), not just the target column, so column-local
# "insert a row"/"delete a row" semantics are implemented here by
# hand: read the column into a flat list, mutate the list, write it
# back, then clear any now-unused trailing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Get-ColumnList($sheet, $colLetter, $startRow, $endRow) {
    $rng = $sheet.Range($colLetter + $startRow + ":" + $colLetter + $endRow)
    $raw = $rng.Value()
    $list = New-Object System.Collections.ArrayList
    $n = $endRow - $startRow + 1
    for ($i = 1; $i -le $n; $i++) {
        [void]$list.Add($raw[$i, 1])
    }
    return $list
}

function Set-ColumnList($sheet, $colLetter, $startRow, $list, $clearThroughRow) {
    $n = $list.Count
    if ($n -gt 0) {
        $arr = New-Object 'object[,]' $n, 1
        for ($i = 0; $i -lt $n; $i++) {
            $v = $list[$i]
            if ($null -eq $v) { $v = "" }
            $arr[$i, 0] = $v
        }
        $endRow = $startRow + $n - 1
        $sheet.Range($colLetter + $startRow + ":" + $colLetter + $endRow).Value2 = $arr
    }
    $clearStart = $startRow + $n
    if ($clearStart -le $clearThroughRow) {
        $sheet.Range($colLetter + $clearStart + ":" + $colLetter + $clearThroughRow).ClearContents()
    }
}

function Find-Index($list, $value) {
    for ($i = 0; $i -lt $list.Count; $i++) {
        if ($list[$i] -eq $value) {
            return $i
        }
    }
    return -1
}

# ---- column F ("base"): drop "clear(variables)" ----------------
$colF = Get-ColumnList $ws "F" 2 40
[void]$colF.Remove("clear(variables)")
Set-ColumnList $ws "F" 2 $colF 40

# ---- column J ("external"): append "terminate(programName)" ----
$colJ = Get-ColumnList $ws "J" 2 6
$colJ[$colJ.Count - 1] = "terminate(programName)"
Set-ColumnList $ws "J" 2 $colJ 6

# ---- column L ("io"): insert "assertPath(path)" -----------------
$colL = Get-ColumnList $ws "L" 2 29
$idx = Find-Index $colL "assertReadableFile(file,minByte)"
$colL.Insert($idx, "assertPath(path)")
Set-ColumnList $ws "L" 2 $colL 30

# ---- column Z ("web"): fix one entry + insert two new ones ------
$colZ = Get-ColumnList $ws "Z" 2 135
$idx4 = Find-Index $colZ "assertAttributeContains(locator,attrName,contains)"
$colZ[$idx4] = "assertAttributeContain(locator,attrName,contains)"
$idxIns = Find-Index $colZ "savePageAsFile(sessionIdName,url,file)"
$colZ.Insert($idxIns + 1, "saveSelectedText(var,locator)")
$colZ.Insert($idxIns + 2, "saveSelectedValue(var,locator)")
Set-ColumnList $ws "Z" 2 $colZ 137

# ---- keep the defined names in sync with the new row counts -----
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("io").RefersTo = "='#system'!`$L`$2:`$L`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
